# Fix the typo in the sheet name "dragonEyire" -> "dragonEyrie"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dragonEyire")
$ws.Name = "dragonEyrie"

# Move the active/selected tab from "toolShop" to "dragonEyrie" and
# update its selected cell to C26 (single cell selection).
$ws.Activate()
$ws.Range("C26").Select()
